# Fixed update to excel issue
# 1) Rename "Requested quantity" header on "Weekly Quantity" -> "Weekly_PO_Qty"
# 2) Rename "Requested quantity" header on "Monthly Trend"   -> "Monthly_PO_Qty"
# 3) Add a new "PO Forecast" sheet (ds / PO_Forecast / yhat_lower / yhat_upper)

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# New "PO Forecast" sheet, added after the last existing sheet so it lands
# as the 3rd tab (sheetId 3 / rId3).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PO Forecast"

# Match the page margins used by the other two sheets (0.75/0.75/1/1/0.5/0.5 in).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"

$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$poForecastData = @(
    @(44948.99999999999, 7, 0.02982871381630952, 14.30086238726629),
    @(44955.99999999999, 7, 0.2340109939040181, 13.92757177276943),
    @(44962.99999999999, 7, -0.2023308452166285, 13.95633725273942),
    @(44976.99999999999, 6, -1.169306809118499, 13.18610882069008),
    @(44983.99999999999, 6, -0.6667802442169024, 13.34106106452969),
    @(44990.99999999999, 6, -1.504018196022507, 12.50029598671055),
    @(44997.99999999999, 5, -2.215664203801743, 12.2417197837185),
    @(45011.99999999999, 5, -2.238481149154647, 12.31023370101259),
    @(45039.99999999999, 4, -3.370834130384118, 10.78664347112487),
    @(45046.99999999999, 4, -3.819881473088953, 10.8001302855283),
    @(45053.99999999999, 3, -3.84251233084203, 10.52831270444468),
    @(45060.99999999999, 3, -3.74751854203385, 10.13941623860719),
    @(45067.99999999999, 3, -4.699670220088018, 9.660264688906278),
    @(45074.99999999999, 2, -4.979157537808404, 9.920579607436947),
    @(45081.99999999999, 2, -4.909593130829765, 9.035364591465026),
    @(45088.99999999999, 2, -5.19445768194213, 9.256696219640475),
    @(45095.99999999999, 2, -5.413110287549483, 8.45014813067527)
)

$r = 2
foreach ($row in $poForecastData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$ws.Range("A2:A18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
